$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the number format (date style) from the last existing data row (A251)
# down into the new date cells A252:A255 so the appended rows match the
# existing column A formatting (style index 2 / custom date format).
$ws.Range("A251").Copy()
$ws.Range("A252:A255").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# New data rows appended to the bottom of the table (rows 252-255).
$data = @(
    @(44326, 3, 27, 105.679282946495),
    @(44327, 2, 29, 113.5073779795687),
    @(44328, 0, 28, 109.5933304630318),
    @(44329, 2, 25, 97.85118791342127)
)

$r = 252
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}
